$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data table: rows 2-7, columns A-T
# A = Sending cluster, B = Ligand symbol, C = Receptor symbol, D = Target cluster
# E..T numeric statistics

$data = @(
    @("ECs","Il18","Il1rapl1","ECs",1,0.3333333333333333,1.926847333333333,5.780542,0.1758475956799574,0.1758475956799574,1,0.3333333333333333,0.02412133333333333,0.072364,0.2229307800605663,0.2229307800605663,0.04647812680977777,0.4183031412879999,0.03920184167670798,0.03920184167670798),
    @("ECs","Il18","Il1rapl1","sCs",1,0.3333333333333333,1.926847333333333,5.780542,0.1758475956799574,0.1758475956799574,3,1,0.08407966666666666,0.252239,0.7770692199394338,0.7770692199394338,0.1620086815042222,1.458078133538,0.1366457540032495,0.1366457540032495),
    @("FAPs","Il18","Il1rapl1","ECs",3,1,6.602185666666666,19.806557,0.6025274839535862,0.6025274839535861,1,0.3333333333333333,0.02412133333333333,0.072364,0.2229307800605663,0.2229307800605663,0.1592535211942222,1.433281690748,0.1343219220057033,0.1343219220057033),
    @("FAPs","Il18","Il1rapl1","sCs",3,1,6.602185666666666,19.806557,0.6025274839535862,0.6025274839535861,3,1,0.08407966666666666,0.252239,0.7770692199394338,0.7770692199394338,0.5551095701247777,4.995986131122999,0.468205561947883,0.4682055619478829),
    @("sCs","Il18","Il1rapl1","ECs",3,1,2.428451666666667,7.285355,0.2216249203664564,0.2216249203664564,1,0.3333333333333333,0.02412133333333333,0.072364,0.2229307800605663,0.2229307800605663,0.05857749213555555,0.52719742922,0.04940701637815501,0.04940701637815501),
    @("sCs","Il18","Il1rapl1","sCs",3,1,2.428451666666667,7.285355,0.2216249203664564,0.2216249203664564,3,1,0.08407966666666666,0.252239,0.7770692199394338,0.7770692199394338,0.2041834066494444,1.837650659845,0.1722179039883014,0.1722179039883014)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
